$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(4, 8).Value = 92.7  # H4: was 110.333336
$ws.Cells.Item(4, 9).Value = 96.14286  # I4: was 110.333336
$ws.Cells.Item(4, 10).Value = 84.666664  # J4: was 0
$ws.Cells.Item(4, 11).Value = 96.14286  # K4: was 110.333336
$ws.Cells.Item(4, 12).Value = 84.666664  # L4: was 0
$ws.Cells.Item(4, 13).Value = 17.85714  # M4: was 3.666663999999997
$ws.Cells.Item(4, 14).Value = -312.666664  # N4: was None
$ws.Cells.Item(18, 8).Value = 500  # H18: was 466.66666
$ws.Cells.Item(18, 9).Value = 500  # I18: was 466.66666
$ws.Cells.Item(18, 11).Value = 500  # K18: was 466.66666
$ws.Cells.Item(18, 13).Value = -216  # M18: was -182.66666
$ws.Cells.Item(33, 8).Value = 202.625  # H33: was 167.95
$ws.Cells.Item(33, 9).Value = 202.625  # I33: was 175.63158
$ws.Cells.Item(33, 10).Value = 0  # J33: was 22
$ws.Cells.Item(33, 11).Value = 202.625  # K33: was 175.63158
$ws.Cells.Item(33, 12).Value = 0  # L33: was 22
$ws.Cells.Item(33, 13).Value = 26.375  # M33: was 53.36841999999999
$ws.Cells.Item(33, 14).ClearContents()  # N33: was -480
$ws.Cells.Item(62, 8).Value = 2685.261  # H62: was 3038.55
$ws.Cells.Item(62, 9).Value = 1969.2307  # I62: was 2655.5557
$ws.Cells.Item(62, 10).Value = 3616.1  # J62: was 3351.9092
$ws.Cells.Item(62, 11).Value = 1969.2307  # K62: was 2655.5557
$ws.Cells.Item(62, 12).Value = 3616.1  # L62: was 3351.9092
$ws.Cells.Item(62, 13).Value = -1345.2307  # M62: was -2031.5557
$ws.Cells.Item(62, 14).Value = -4864.1  # N62: was -4599.9092
$ws.Cells.Item(65, 8).Value = 2685.261  # H65: was 3038.55
$ws.Cells.Item(65, 9).Value = 1969.2307  # I65: was 2655.5557
$ws.Cells.Item(65, 10).Value = 3616.1  # J65: was 3351.9092
$ws.Cells.Item(65, 11).Value = 9846.1535  # K65: was 13277.7785
$ws.Cells.Item(65, 12).Value = 18080.5  # L65: was 16759.546
$ws.Cells.Item(65, 13).Value = -6726.1535  # M65: was -10157.7785
$ws.Cells.Item(65, 14).Value = -24320.5  # N65: was -22999.546
$ws.Cells.Item(98, 8).Value = 737  # H98: was 780.90625
$ws.Cells.Item(98, 9).Value = 737  # I98: was 634.96155
$ws.Cells.Item(98, 10).Value = 0  # J98: was 1413.3334
$ws.Cells.Item(98, 11).Value = 737  # K98: was 634.96155
$ws.Cells.Item(98, 12).Value = 0  # L98: was 1413.3334
$ws.Cells.Item(98, 13).Value = 761  # M98: was 863.03845
$ws.Cells.Item(98, 14).ClearContents()  # N98: was -4409.3334
$ws.Cells.Item(122, 8).Value = 737  # H122: was 780.90625
$ws.Cells.Item(122, 9).Value = 737  # I122: was 634.96155
$ws.Cells.Item(122, 10).Value = 0  # J122: was 1413.3334
$ws.Cells.Item(122, 11).Value = 2211  # K122: was 1904.88465
$ws.Cells.Item(122, 12).Value = 0  # L122: was 4240.0002
$ws.Cells.Item(122, 13).Value = 239  # M122: was 545.11535
$ws.Cells.Item(122, 14).ClearContents()  # N122: was -9140.0002
$ws.Cells.Item(129, 8).Value = 278482.47  # H129: was 257121.89
$ws.Cells.Item(129, 10).Value = 334093.97  # J129: was 303794.06
$ws.Cells.Item(129, 12).Value = 1002281.91  # L129: was 911382.1799999999
$ws.Cells.Item(129, 14).Value = -1012281.91  # N129: was -921382.1799999999
$ws.Cells.Item(137, 8).Value = 1715.909  # H137: was 1567.6428
$ws.Cells.Item(137, 9).Value = 1820  # I137: was 1612.2
$ws.Cells.Item(137, 10).Value = 1492.8572  # J137: was 1456.25
$ws.Cells.Item(137, 11).Value = 5460  # K137: was 4836.6
$ws.Cells.Item(137, 12).Value = 4478.571599999999  # L137: was 4368.75
$ws.Cells.Item(137, 13).Value = -2910  # M137: was -2286.6
$ws.Cells.Item(137, 14).Value = -9578.571599999999  # N137: was -9468.75
$ws.Cells.Item(138, 8).Value = 2505.919  # H138: was 2445.2368
$ws.Cells.Item(138, 9).Value = 1462.7778  # I138: was 1396.3158
$ws.Cells.Item(138, 11).Value = 4388.3334  # K138: was 4188.9474
$ws.Cells.Item(138, 13).Value = 751.6665999999996  # M138: was 951.0526

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 3396.6296  # H45: was 3661.6538
$ws.Cells.Item(45, 9).Value = 3508.0908  # I45: was 3739.7
$ws.Cells.Item(45, 10).Value = 3320  # J45: was 3612.875
$ws.Cells.Item(45, 11).Value = 3508.0908  # K45: was 3739.7
$ws.Cells.Item(45, 12).Value = 3320  # L45: was 3612.875
$ws.Cells.Item(45, 13).Value = -3131.0908  # M45: was -3362.7
$ws.Cells.Item(45, 14).Value = -4074  # N45: was -4366.875
$ws.Cells.Item(61, 8).Value = 3199.4517  # H61: was 3128.2258
$ws.Cells.Item(61, 9).Value = 2983.32  # I61: was 2891.3462
$ws.Cells.Item(61, 10).Value = 4100  # J61: was 4360
$ws.Cells.Item(61, 11).Value = 2983.32  # K61: was 2891.3462
$ws.Cells.Item(61, 12).Value = 4100  # L61: was 4360
$ws.Cells.Item(61, 13).Value = -2771.32  # M61: was -2679.3462
$ws.Cells.Item(61, 14).Value = -4524  # N61: was -4784
$ws.Cells.Item(88, 8).Value = 126794.125  # H88: was 201810.8
$ws.Cells.Item(88, 9).Value = 2125  # I88: was 2333.3333
$ws.Cells.Item(88, 10).Value = 251463.25  # J88: was 501027
$ws.Cells.Item(88, 11).Value = 2125  # K88: was 2333.3333
$ws.Cells.Item(88, 12).Value = 251463.25  # L88: was 501027
$ws.Cells.Item(88, 13).Value = -1719  # M88: was -1927.3333
$ws.Cells.Item(88, 14).Value = -252275.25  # N88: was -501839
$ws.Cells.Item(91, 8).Value = 126794.125  # H91: was 201810.8
$ws.Cells.Item(91, 9).Value = 2125  # I91: was 2333.3333
$ws.Cells.Item(91, 10).Value = 251463.25  # J91: was 501027
$ws.Cells.Item(91, 11).Value = 2125  # K91: was 2333.3333
$ws.Cells.Item(91, 12).Value = 251463.25  # L91: was 501027
$ws.Cells.Item(91, 13).Value = -721  # M91: was -929.3332999999998
$ws.Cells.Item(91, 14).Value = -254271.25  # N91: was -503835
$ws.Cells.Item(122, 8).Value = 1633.5883  # H122: was 1740.1333
$ws.Cells.Item(122, 9).Value = 1718.0667  # I122: was 1793
$ws.Cells.Item(122, 11).Value = 5154.2001  # K122: was 5379
$ws.Cells.Item(122, 13).Value = -2704.2001  # M122: was -2929
$ws.Cells.Item(132, 8).Value = 16189.527  # H132: was 17108.97
$ws.Cells.Item(132, 9).Value = 2157.087  # I132: was 2365.25
$ws.Cells.Item(132, 10).Value = 41016.152  # J132: was 38171.43
$ws.Cells.Item(132, 11).Value = 6471.261  # K132: was 7095.75
$ws.Cells.Item(132, 12).Value = 123048.456  # L132: was 114514.29
$ws.Cells.Item(132, 13).Value = -3941.261  # M132: was -4565.75
$ws.Cells.Item(132, 14).Value = -128108.456  # N132: was -119574.29
$ws.Cells.Item(136, 8).Value = 3199.4517  # H136: was 3128.2258
$ws.Cells.Item(136, 9).Value = 2983.32  # I136: was 2891.3462
$ws.Cells.Item(136, 10).Value = 4100  # J136: was 4360
$ws.Cells.Item(136, 11).Value = 8949.960000000001  # K136: was 8674.0386
$ws.Cells.Item(136, 12).Value = 12300  # L136: was 13080
$ws.Cells.Item(136, 13).Value = -6399.960000000001  # M136: was -6124.0386
$ws.Cells.Item(136, 14).Value = -17400  # N136: was -18180

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 1725.4103  # H86: was 1776.7142
$ws.Cells.Item(86, 9).Value = 1571.64  # I86: was 1608.409
$ws.Cells.Item(86, 10).Value = 2000  # J86: was 2061.5386
$ws.Cells.Item(86, 11).Value = 1571.64  # K86: was 1608.409
$ws.Cells.Item(86, 12).Value = 2000  # L86: was 2061.5386
$ws.Cells.Item(86, 13).Value = -448.6400000000001  # M86: was -485.4090000000001
$ws.Cells.Item(86, 14).Value = -4246  # N86: was -4307.5386
$ws.Cells.Item(89, 8).Value = 1725.4103  # H89: was 1776.7142
$ws.Cells.Item(89, 9).Value = 1571.64  # I89: was 1608.409
$ws.Cells.Item(89, 10).Value = 2000  # J89: was 2061.5386
$ws.Cells.Item(89, 11).Value = 7858.200000000001  # K89: was 8042.045
$ws.Cells.Item(89, 12).Value = 10000  # L89: was 10307.693
$ws.Cells.Item(89, 13).Value = -2242.200000000001  # M89: was -2426.045
$ws.Cells.Item(89, 14).Value = -21232  # N89: was -21539.693
$ws.Cells.Item(94, 8).Value = 596.4194  # H94: was 673.3333
$ws.Cells.Item(94, 9).Value = 424.73914  # I94: was 497.89474
$ws.Cells.Item(94, 11).Value = 424.73914  # K94: was 497.89474
$ws.Cells.Item(94, 13).Value = 26.26085999999998  # M94: was -46.89474000000001
$ws.Cells.Item(134, 8).Value = 3155.3408  # H134: was 3322.5557
$ws.Cells.Item(134, 9).Value = 3083.425  # I134: was 3085.425
$ws.Cells.Item(134, 10).Value = 3874.5  # J134: was 5219.6
$ws.Cells.Item(134, 11).Value = 9250.275000000001  # K134: was 9256.275000000001
$ws.Cells.Item(134, 12).Value = 11623.5  # L134: was 15658.8
$ws.Cells.Item(134, 13).Value = -6715.275000000001  # M134: was -6721.275000000001
$ws.Cells.Item(134, 14).Value = -16693.5  # N134: was -20728.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4250.4814  # H31: was 3751.7097
$ws.Cells.Item(31, 9).Value = 996.25  # I31: was 939.2308
$ws.Cells.Item(31, 10).Value = 5620.684  # J31: was 5782.9443
$ws.Cells.Item(31, 11).Value = 996.25  # K31: was 939.2308
$ws.Cells.Item(31, 12).Value = 5620.684  # L31: was 5782.9443
$ws.Cells.Item(31, 13).Value = -701.25  # M31: was -644.2308
$ws.Cells.Item(31, 14).Value = -6210.684  # N31: was -6372.9443
$ws.Cells.Item(34, 8).Value = 4250.4814  # H34: was 3751.7097
$ws.Cells.Item(34, 9).Value = 996.25  # I34: was 939.2308
$ws.Cells.Item(34, 10).Value = 5620.684  # J34: was 5782.9443
$ws.Cells.Item(34, 11).Value = 996.25  # K34: was 939.2308
$ws.Cells.Item(34, 12).Value = 5620.684  # L34: was 5782.9443
$ws.Cells.Item(34, 13).Value = -794.25  # M34: was -737.2308
$ws.Cells.Item(34, 14).Value = -6024.684  # N34: was -6186.9443
$ws.Cells.Item(105, 8).Value = 1144.5883  # H105: was 1189.3684
$ws.Cells.Item(105, 9).Value = 974.0714  # I105: was 969.13336
$ws.Cells.Item(105, 10).Value = 1940.3334  # J105: was 2015.25
$ws.Cells.Item(105, 11).Value = 974.0714  # K105: was 969.13336
$ws.Cells.Item(105, 12).Value = 1940.3334  # L105: was 2015.25
$ws.Cells.Item(105, 13).Value = 772.9286  # M105: was 777.86664
$ws.Cells.Item(105, 14).Value = -5434.3334  # N105: was -5509.25
$ws.Cells.Item(132, 8).Value = 2387.6428  # H132: was 2645.28
$ws.Cells.Item(132, 9).Value = 1402.1052  # I132: was 1528.8889
$ws.Cells.Item(132, 10).Value = 4468.222  # J132: was 5516
$ws.Cells.Item(132, 11).Value = 4206.3156  # K132: was 4586.6667
$ws.Cells.Item(132, 12).Value = 13404.666  # L132: was 16548
$ws.Cells.Item(132, 13).Value = -1676.3156  # M132: was -2056.6667
$ws.Cells.Item(132, 14).Value = -18464.666  # N132: was -21608

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 730.0599999999999  # H131: was 736.0303
$ws.Cells.Item(131, 10).Value = 734.80414  # J131: was 741.01044
$ws.Cells.Item(131, 12).Value = 2204.41242  # L131: was 2223.03132
$ws.Cells.Item(131, 14).Value = -12284.41242  # N131: was -12303.03132

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 1002.2  # H97: was 1344.7
$ws.Cells.Item(97, 9).Value = 1023.2  # I97: was 1288.5555
$ws.Cells.Item(97, 10).Value = 918.2  # J97: was 1850
$ws.Cells.Item(97, 11).Value = 1023.2  # K97: was 1288.5555
$ws.Cells.Item(97, 12).Value = 918.2  # L97: was 1850
$ws.Cells.Item(97, 13).Value = -527.2  # M97: was -792.5554999999999
$ws.Cells.Item(97, 14).Value = -1910.2  # N97: was -2842
$ws.Cells.Item(102, 8).Value = 1603.4138  # H102: was 1762.6086
$ws.Cells.Item(102, 9).Value = 1323.5217  # I102: was 1396.9474
$ws.Cells.Item(102, 10).Value = 2676.3333  # J102: was 3499.5
$ws.Cells.Item(102, 11).Value = 1323.5217  # K102: was 1396.9474
$ws.Cells.Item(102, 12).Value = 2676.3333  # L102: was 3499.5
$ws.Cells.Item(102, 13).Value = 298.4783  # M102: was 225.0526
$ws.Cells.Item(102, 14).Value = -5920.3333  # N102: was -6743.5
$ws.Cells.Item(122, 8).Value = 2527.0667  # H122: was 2432
$ws.Cells.Item(122, 9).Value = 1113.375  # I122: was 1101.4445
$ws.Cells.Item(122, 11).Value = 3340.125  # K122: was 3304.3335
$ws.Cells.Item(122, 13).Value = -890.125  # M122: was -854.3335000000002
$ws.Cells.Item(132, 8).Value = 32918.5  # H132: was 41845.215
$ws.Cells.Item(132, 9).Value = 6174.5454  # I132: was 8745.714
$ws.Cells.Item(132, 11).Value = 18523.6362  # K132: was 26237.142
$ws.Cells.Item(132, 13).Value = -15993.6362  # M132: was -23707.142

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3814.6667  # H7: was 4010.353
$ws.Cells.Item(7, 9).Value = 3557.3333  # I7: was 4120
$ws.Cells.Item(7, 10).Value = 4329.3335  # J7: was 3853.7144
$ws.Cells.Item(7, 11).Value = 3557.3333  # K7: was 4120
$ws.Cells.Item(7, 12).Value = 4329.3335  # L7: was 3853.7144
$ws.Cells.Item(7, 13).Value = -3445.3333  # M7: was -4008
$ws.Cells.Item(7, 14).Value = -4553.3335  # N7: was -4077.7144
$ws.Cells.Item(22, 8).Value = 2842.45  # H22: was 2588.5908
$ws.Cells.Item(22, 9).Value = 3920.9167  # I22: was 3034.4375
$ws.Cells.Item(22, 10).Value = 1224.75  # J22: was 1399.6666
$ws.Cells.Item(22, 11).Value = 3920.9167  # K22: was 3034.4375
$ws.Cells.Item(22, 12).Value = 1224.75  # L22: was 1399.6666
$ws.Cells.Item(22, 13).Value = -3625.9167  # M22: was -2739.4375
$ws.Cells.Item(22, 14).Value = -1814.75  # N22: was -1989.6666
$ws.Cells.Item(27, 8).Value = 2842.45  # H27: was 2588.5908
$ws.Cells.Item(27, 9).Value = 3920.9167  # I27: was 3034.4375
$ws.Cells.Item(27, 10).Value = 1224.75  # J27: was 1399.6666
$ws.Cells.Item(27, 11).Value = 3920.9167  # K27: was 3034.4375
$ws.Cells.Item(27, 12).Value = 1224.75  # L27: was 1399.6666
$ws.Cells.Item(27, 13).Value = -3813.9167  # M27: was -2927.4375
$ws.Cells.Item(27, 14).Value = -1438.75  # N27: was -1613.6666
$ws.Cells.Item(122, 8).Value = 787251.0600000001  # H122: was 787270.7
$ws.Cells.Item(122, 9).Value = 936479.9  # I122: was 936503.2
$ws.Cells.Item(122, 11).Value = 2809439.7  # K122: was 2809509.6
$ws.Cells.Item(122, 13).Value = -2806989.7  # M122: was -2807059.6
$ws.Cells.Item(126, 8).Value = 3814.6667  # H126: was 4010.353
$ws.Cells.Item(126, 9).Value = 3557.3333  # I126: was 4120
$ws.Cells.Item(126, 10).Value = 4329.3335  # J126: was 3853.7144
$ws.Cells.Item(126, 11).Value = 10671.9999  # K126: was 12360
$ws.Cells.Item(126, 12).Value = 12988.0005  # L126: was 11561.1432
$ws.Cells.Item(126, 13).Value = -8201.999899999999  # M126: was -9890
$ws.Cells.Item(126, 14).Value = -17928.0005  # N126: was -16501.1432
$ws.Cells.Item(136, 8).Value = 1335.4445  # H136: was 1392.04
$ws.Cells.Item(136, 9).Value = 1335.4445  # I136: was 1392.04
$ws.Cells.Item(136, 11).Value = 4006.3335  # K136: was 4176.12
$ws.Cells.Item(136, 13).Value = -1456.3335  # M136: was -1626.12

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 1394.1936  # H132: was 1514.8214
$ws.Cells.Item(132, 9).Value = 1185.4117  # I132: was 1310.1333
$ws.Cells.Item(132, 10).Value = 1647.7142  # J132: was 1751
$ws.Cells.Item(132, 11).Value = 3556.2351  # K132: was 3930.3999
$ws.Cells.Item(132, 12).Value = 4943.142599999999  # L132: was 5253
$ws.Cells.Item(132, 13).Value = -1026.2351  # M132: was -1400.3999
$ws.Cells.Item(132, 14).Value = -10003.1426  # N132: was -10313
$ws.Cells.Item(136, 8).Value = 18870132  # H136: was 16951084
$ws.Cells.Item(136, 9).Value = 23256684  # I136: was 20408944
$ws.Cells.Item(136, 10).Value = 7951.5  # J136: was 7572
$ws.Cells.Item(136, 11).Value = 69770052  # K136: was 61226832
$ws.Cells.Item(136, 12).Value = 23854.5  # L136: was 22716
$ws.Cells.Item(136, 13).Value = -69767502  # M136: was -61224282
$ws.Cells.Item(136, 14).Value = -28954.5  # N136: was -27816
$ws.Cells.Item(140, 8).Value = 44719.8  # H140: was 45199.75
$ws.Cells.Item(140, 10).Value = 44719.8  # J140: was 45199.75
$ws.Cells.Item(140, 12).Value = 44719.8  # L140: was 45199.75
$ws.Cells.Item(140, 14).Value = -55079.8  # N140: was -55559.75
$ws.Cells.Item(141, 8).Value = 71928.75  # H141: was 72197.8
$ws.Cells.Item(141, 10).Value = 71928.75  # J141: was 72197.8
$ws.Cells.Item(141, 12).Value = 71928.75  # L141: was 72197.8
$ws.Cells.Item(141, 14).Value = -82288.75  # N141: was -82557.8
